# Apply the "roolit.xlsx" permission-matrix update (YIT+Destia+Lemminkäinen
# proposal) to the "Oikeudet" sheet.
#
# Summary of content changes (see xl/revisions/revisionLog18-20.xml in the
# target OOXML for the authoritative cell-by-cell trail):
#   - V6   : "Urakan turvallisuusvastaava" -> "Kelikeskus"  (role renamed)
#   - S7   : "R" -> "R,W"
#   - V7   : "R" -> "R,W"
#   - S13  : ""  -> "R,W"   (also gains the thin-border / centered / automatic
#   - S14  : ""  -> "R,W"    font-colour style already used by the rest of
#   - S15  : ""  -> "R,W"    that row's permission columns)
#   - S16  : "R" -> "R,W"
#   - S17  : "R" -> "R,W"
#   - S18  : ""  -> "R,W"
#   - S19  : "R" -> "R,W"
#   - S22  : ""  -> "R"
#   - V22  : "R" -> ""  (cleared)
#   - V23  : "R" -> ""  (cleared)
#   - S25  : ""  -> "R,W"  (also restyled to match the rest of the column)
#   - V25  : "R,W" -> ""  (cleared)
#   - S46  : ""  -> "R"
#   - S48  : ""  -> "R"
#   - S53  : ""  -> "R"
#   - S63  : ""  -> "R,W"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Oikeudet")

# Role renamed: "Urakan turvallisuusvastaava" -> "Kelikeskus"
$ws.Range("V6").Value = "Kelikeskus"

# Simple value-only updates (existing cell formatting is left untouched).
$ws.Range("S7").Value  = "R,W"
$ws.Range("V7").Value  = "R,W"
$ws.Range("S16").Value = "R,W"
$ws.Range("S17").Value = "R,W"
$ws.Range("S19").Value = "R,W"
$ws.Range("S22").Value = "R"
$ws.Range("S46").Value = "R"
$ws.Range("S48").Value = "R"
$ws.Range("S53").Value = "R"
$ws.Range("S63").Value = "R,W"

# Cells that were blank and now need both a new value and the bordered /
# centered / automatic-colour look already used by the rest of their row
# (matches the R column in the same row). Cloning the format from the
# neighbouring cell is the most reliable way to reproduce that exactly.
$newlyFilled = @{
    "S13" = "R13"
    "S14" = "R14"
    "S15" = "R15"
    "S18" = "R18"
    "S25" = "R25"
}
foreach ($target in $newlyFilled.Keys) {
    $src = $ws.Range($newlyFilled[$target])
    $dst = $ws.Range($target)
    $src.Copy()
    $dst.PasteSpecial(-4122)  # xlPasteFormats
    $dst.Value = "R,W"
}
$excel.CutCopyMode = 0

# Cells whose value is removed (role no longer applies there); keep formatting.
$ws.Range("V22").ClearContents()
$ws.Range("V23").ClearContents()
$ws.Range("V25").ClearContents()
